$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Move the data block from column E down to column G (rows 2-7).
# ------------------------------------------------------------------
$ws.Range("E2:E7").Cut($ws.Range("G2:G7"))

# ------------------------------------------------------------------
# 2) Re-font everything to Arial 10 (was Aptos Narrow 11).
#    A1 keeps its bold weight, everything else stays regular.
#    Temporarily un-bold A1 so that it shares the same "regular"
#    font-derivation path as the rest of the sheet - this lets the
#    engine re-use a single interned font object instead of forking
#    a whole separate bold-track font lineage.
# ------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $false

$ws.Range("A1:A2").Font.Name = "Arial"
$ws.Range("A1:A2").Font.Size = 10

$ws.Range("G2:G7").Font.Name = "Arial"
$ws.Range("G2:G7").Font.Size = 10

# Now restore the bold weight on the header cell (derives a single new
# "Arial 10 bold" font from the already-built "Arial 10" font).
$ws.Range("A1").Font.Bold = $true

# ------------------------------------------------------------------
# 3) Restore the on-screen selection to E3 (matches the saved view).
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("E3").Select()
